$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 9088
$ws1.Range("F4").Value = 469
$ws1.Range("F5").Value = 454

# Sheet "全部类型" (All types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 9088
$ws4.Range("F4").Value = 469
$ws4.Range("F6").Value = 454
